# [plugin-excel] Fix loading of non-string cells by FROM_EXCEL table transformer
#
# Adds a new "DifferentTypes" worksheet (after "RepeatingData") that exercises
# boolean / numeric / string / formula cell types, used by the FROM_EXCEL
# table-transformer regression test. Also tidies up the duplicated "_xlcn."
# worksheet-connection defined names that Excel had suffixed with a stray
# trailing "1".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Strip the stray trailing "1" Excel had appended to these auto-generated
#    Power Pivot worksheet-connection defined names.
# ---------------------------------------------------------------------------
foreach ($dn in $wb.Names) {
    if ($dn.Name.EndsWith("1") -and $dn.Name.StartsWith("_xlcn.WorksheetConnection_")) {
        $dn.Name = $dn.Name.Substring(0, $dn.Name.Length - 1)
    }
}

# ---------------------------------------------------------------------------
# 2. Add the new "DifferentTypes" worksheet as the last (rightmost) tab.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "DifferentTypes"

# Header row.
$ws.Range("A1").Value = "Boolean"
$ws.Range("B1").Value = "Number"
$ws.Range("C1").Value = "String"
$ws.Range("D1").Value = "Formula"

# Row 2.
$ws.Range("A2").Value = $true
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "STRING"
$ws.Range("D2").Formula = "=B2+B3"

# Row 3.
$ws.Range("A3").Value = $false
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "string"
$ws.Range("D3").Formula = "=C2+C3"

# Make the new sheet the active tab/selection, mirroring the recorded
# template (A1:D3 selected on the newly active "DifferentTypes" sheet).
$ws.Activate() | Out-Null
$ws.Range("A1:D3").Select() | Out-Null

$wb.Save()
